{"js": "// Apply the \"Euphoria Megaways\" copy-rewrite described by the commit\n// \"Added many more features\": retitle the page, swap the \"What we like\" /\n// \"What we don't like\" bullet copy, and refresh the closing summary\n// (bold title + italic teaser) to match.\n//\n// Several of the old phrases (e.g. \"117,649 ways to win\") also occur as a\n// substring inside unrelated narrative sentences elsewhere in the body, so\n// a document-wide Body.search() would over-match. Instead:\n//   1. find the paragraph(s) whose *entire* text equals the old string,\n//   2. search() *within* that paragraph for the old text to get a Range\n//      scoped to just the run(s) holding it, then\n//   3. Range.insertText(..., Word.InsertLocation.replace) to swap only\n//      that text.\n// This leaves paragraph styles/formatting and any neighbouring empty runs\n// (e.g. the leading <w:r/>) untouched, matching the diff exactly.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Ordered old -> new text pairs, exactly as in the diff. Order matters:\n// \"High volatility\" is itself introduced as *new* text by one\n// replacement, and is the *old* text consumed by a later one, so we must\n// process them top-to-bottom and re-read paragraph text between steps.\nconst replacements = [\n    [\n        \"Play Euphoria Megaways for Free: Exciting Bonus Games\",\n        \"Play Free: Euphoria Megaways Slot Review\"\n    ],\n    [\"117,649 ways to win\", \"Variety of special symbols and features\"],\n    [\"Exciting bonus games\", \"High winning potential\"],\n    [\"Free Spins and Super Free Spins\", \"Smooth gameplay on mobile devices\"],\n    [\n        \"High maximum win potential\",\n        \"Ability to play for free on the provider's website\"\n    ],\n    [\"High volatility\", \"Limited number of Scatter symbols\"],\n    [\"Autoplay function requires manual limits\", \"High volatility\"],\n    [\n        \"Read our review of Euphoria Megaways, the high-volatility slot game with free spins and exciting bonus games. Try it for free on the provider's website.\",\n        \"Read our review of Euphoria Megaways slot game. Play for free and discover its features.\"\n    ]\n];\n\nfor (const [oldText, newText] of replacements) {\n    for (let i = 0; i < paragraphs.items.length; i++) {\n        const paragraph = paragraphs.items[i];\n        // Strip the trailing paragraph-mark CR that Word.js includes in .text.\n        const currentText = paragraph.text.replace(/\\r$/, \"\");\n        if (currentText === oldText) {\n            const hits = paragraph.search(oldText, { matchCase: true });\n            hits.load(\"text\");\n            await context.sync();\n            for (let j = 0; j < hits.items.length; j++) {\n                hits.items[j].insertText(newText, Word.InsertLocation.replace);\n            }\n        }\n    }\n    // Sync + reload so the next replacement (some old strings, like\n    // \"High volatility\", are themselves produced by an earlier step) sees\n    // up-to-date paragraph text.\n    await context.sync();\n    paragraphs.load(\"text\");\n    await context.sync();\n}\n", "ps1": "# Apply the \"Euphoria Megaways\" copy-rewrite described by the commit\n# \"Added many more features\": retitle the page, swap the \"What we like\" /\n# \"What we don't like\" bullet copy, and refresh the closing summary\n# (bold title + italic teaser) to match.\n#\n# Several of the old phrases (e.g. \"117,649 ways to win\") also occur as a\n# substring inside unrelated narrative sentences elsewhere in the body, so\n# a document-wide Find/Replace over $d.Content would over-match. Instead,\n# for each old/new pair we find the paragraph(s) whose *entire* text\n# equals the old string and run Find/Replace scoped to just that\n# paragraph's Range.\n\n$d = $word.ActiveDocument\n\n# Ordered old -> new text pairs, exactly as in the diff. Order matters:\n# \"High volatility\" is itself introduced as *new* text by one\n# replacement, and is the *old* text consumed by a later one, so this\n# list is processed top-to-bottom against the live document.\n$replacements = @(\n    @(\"Play Euphoria Megaways for Free: Exciting Bonus Games\", \"Play Free: Euphoria Megaways Slot Review\"),\n    @(\"117,649 ways to win\", \"Variety of special symbols and features\"),\n    @(\"Exciting bonus games\", \"High winning potential\"),\n    @(\"Free Spins and Super Free Spins\", \"Smooth gameplay on mobile devices\"),\n    @(\"High maximum win potential\", \"Ability to play for free on the provider's website\"),\n    @(\"High volatility\", \"Limited number of Scatter symbols\"),\n    @(\"Autoplay function requires manual limits\", \"High volatility\"),\n    @(\"Read our review of Euphoria Megaways, the high-volatility slot game with free spins and exciting bonus games. Try it for free on the provider's website.\", \"Read our review of Euphoria Megaways slot game. Play for free and discover its features.\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    foreach ($p in $d.Paragraphs) {\n        # Trim the trailing paragraph-mark (CR) character before comparing.\n        $currentText = $p.Range.Text.TrimEnd([char]13)\n        if ($currentText -eq $oldText) {\n            $rng = $p.Range\n            $rng.Find.ClearFormatting()\n            $rng.Find.Replacement.ClearFormatting()\n            $rng.Find.Text = $oldText\n            $rng.Find.Replacement.Text = $newText\n            $rng.Find.Execute([ref]$null, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$null, [ref]2) | Out-Null\n        }\n    }\n}\n"}
